$d = $word.ActiveDocument
$t = $d.Tables(1)
$row = $t.Rows(3)

# Column 1: date/time stamp
$row.Cells(1).Range.Text = "Nov 20, 10:49pm"

# Column 2: version
$row.Cells(2).Range.Text = "2"

# Column 3: test standard
$row.Cells(3).Range.Text = "Moisture accuracy test in " + [char]0x201C + "dry" + [char]0x201D + " and " + [char]0x201C + "wet" + [char]0x201D + " soil (refer to design document)"

# Column 4: results - first line of text, then two more text paragraphs
# separated by blank paragraphs, mirroring the row above it.
$resultsCell = $row.Cells(4)
$resultsCell.Range.Text = [char]0x201C + "Dry" + [char]0x201D + " soil: average VWC 7%"
$resultsCell.Range.InsertParagraphAfter()
$resultsCell.Range.InsertParagraphAfter()
$resultsCell.Range.InsertParagraphAfter()
$resultsCell.Range.InsertParagraphAfter()
$resultsCell.Range.Paragraphs(3).Range.Text = [char]0x201C + "Wet" + [char]0x201D + " soil: average VWC 93%"
$resultsCell.Range.Paragraphs(5).Range.Text = "Conclusion: passed"
